# "food list in assets" — remove the empty "watermelon" row (row 48, which
# only had a Glycemic-Index label and no Serving Size / Carb data) so every
# remaining food row (chickpeas .. pancakes) shifts up by one and the sheet
# ends at row 62 instead of row 63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the whole row shifts rows 49:63 up to 48:62 automatically,
# carrying their existing values (A/B/C/D) along with them — matching the
# diff exactly (row 48 becomes the old row 49's data, etc.) and shrinking
# the used range from A1:E63 to A1:E62.
$ws.Rows("48").Delete()

# Match the refreshed selection/viewport from the diff (the author had
# scrolled down to/selected F50 after the edit).
$ws.Range("F50").Select()
